$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.006971674199201135, 0.414613479980299, 0.4465464471237033, 0.6682413090521292, 0.6748539099243737, 51)
    3  = @(0.4355859266737307, 0.699453688851532, 1.429667568063306, 1.195687069455594, 1.124827636930724, 50)
    4  = @(0.5940129335889834, 1.050797193125951, 4.135615467765216, 2.033621269500596, 1.965088213786561, 49)
    5  = @(0.5024337226471743, 1.099887659833765, 4.742322585733941, 2.177687439862282, 2.141357598426261, 48)
    6  = @(0.4027908244544923, 0.9781117343447195, 4.417744944668907, 2.101843225521092, 2.08518955102823, 47)
    7  = @(0.3740519112273883, 0.9832035832409267, 4.94152750655546, 2.222954679375057, 2.220672418386649, 38)
    8  = @(0.3686574517514005, 0.9897903167014467, 5.048307810960813, 2.246843966758887, 2.246965739777832, 37)
    9  = @(0.2563981495475554, 1.371782596377584, 8.685257687481315, 2.947076125158852, 3.012171393108808, 20)
    10 = @(-0.4358536319779404, 1.061643998448061, 5.46603034801007, 2.337954308366626, 2.390760086135826, 13)
    11 = @(0.5199476580913862, 0.536199036767537, 0.3603403617065257, 0.600283567746549, 0.3354004967047021, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
